# Updates the division-problem worksheet numbers in the table cells.
# Find.Execute signature used below:
#   Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#           MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#           Format, ReplaceWith, Replace)
# Wrap=1 (wdFindContinue), Replace=2 (wdReplaceAll)

$d = $word.ActiveDocument

$d.Content.Find.Execute("396÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "774÷8=", 2)
$d.Content.Find.Execute("595÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "959÷3=", 2)
$d.Content.Find.Execute("883÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "700÷8=", 2)
$d.Content.Find.Execute("764÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "704÷8=", 2)
$d.Content.Find.Execute("245÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "200÷5=", 2)
$d.Content.Find.Execute("896÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "606÷3=", 2)
$d.Content.Find.Execute("813÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "120÷4=", 2)
$d.Content.Find.Execute("298÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "572÷3=", 2)
$d.Content.Find.Execute("712÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "265÷7=", 2)
$d.Content.Find.Execute("246÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "654÷6=", 2)
$d.Content.Find.Execute("285÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "292÷2=", 2)
$d.Content.Find.Execute("790÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "254÷6=", 2)
$d.Content.Find.Execute("622÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "368÷2=", 2)
$d.Content.Find.Execute("969÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "283÷8=", 2)
$d.Content.Find.Execute("730÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "703÷3=", 2)
$d.Content.Find.Execute("809÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "930÷8=", 2)
$d.Content.Find.Execute("414÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "408÷9=", 2)
$d.Content.Find.Execute("972÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "677÷8=", 2)
$d.Content.Find.Execute("366÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "425÷3=", 2)
$d.Content.Find.Execute("740÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "106÷3=", 2)
$d.Content.Find.Execute("240÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "521÷9=", 2)
$d.Content.Find.Execute("764÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "764÷5=", 2)
$d.Content.Find.Execute("408÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "174÷7=", 2)
$d.Content.Find.Execute("913÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "337÷9=", 2)
$d.Content.Find.Execute("759÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "392÷6=", 2)
